# Refresh the "cryptos" price/volume table with the latest scrape.
# (GitHub Actions bot commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.956.38'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.777.77'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.25'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.550'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.18'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.288'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0703'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0937'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('D12').Value = '2.032.78'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '1.779.33'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.91'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.14%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.620'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '33.926.42'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.14'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.74'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.87'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '0.0₃0783'
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.69'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.08'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.52'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.30'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.08'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.24'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.65'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.51'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '1.393.60'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.655'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.36%  '
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.24'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.35'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.910'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '77.67'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.32'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +13.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.07'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('D46').Value = '0.0₆0137'
$ws.Range('E46').Value = '  +12.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0496'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '107.88'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').Value = '1.932.44'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('E51').Value = '  +0.52%  '
